$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Library_Formula")

# New indicator rows (CR alzo zero e xra corporate), appended after the
# existing last row (52) following the same pattern as the rows above.
$newIndicators = @("INDICATOR_55", "INDICATOR_56", "INDICATOR_57", "INDICATOR_58", "INDICATOR_60", "INDICATOR_61")

$startRow = 53
for ($i = 0; $i -lt $newIndicators.Count; $i++) {
    $r = $startRow + $i
    $indicator = $newIndicators[$i]

    $cellA = $ws.Range("A" + $r)
    $cellB = $ws.Range("B" + $r)
    $cellC = $ws.Range("C" + $r)
    $cellE = $ws.Range("E" + $r)
    $cellF = $ws.Range("F" + $r)

    # Match the look & feel (font) of the surrounding data rows.
    $cellA.Font.Name = "Trebuchet MS"
    $cellA.Font.Size = 10
    $cellB.Font.Name = "Trebuchet MS"
    $cellB.Font.Size = 10
    $cellC.Font.Name = "Trebuchet MS"
    $cellC.Font.Size = 10
    $cellE.Font.Name = "Trebuchet MS"
    $cellE.Font.Size = 10
    $cellF.Font.Name = "Trebuchet MS"
    $cellF.Font.Size = 10

    $cellA.Value = "CREATE/MODIFY"
    $cellB.Value = "LIB_EWS_IT"
    $cellC.Value = $indicator
    $cellE.Value = "String"
    $cellF.Value = "String"
}

# Bring the "Library_Formula" sheet to the front (becomes the active tab)
# and leave the selection on the last block of newly added cells.
$ws.Activate()
$ws.Range("F56:F58").Select()
